$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "GRT-1"
$ws.Range("B2").Value = 131900.7940583276
$ws.Range("C2").Value = 233.8846893310547
$ws.Range("D2").Value = 233.8776092529297
$ws.Range("E2").Value = 235.66824340820312
$ws.Range("F2").Value = 52.003822326660156
$ws.Range("G2").Value = 57.20133972167969
$ws.Range("H2").Value = 52.23653793334961
$ws.Range("I2").Value = 10028.87890625
$ws.Range("J2").Value = 10816.6875
$ws.Range("K2").Value = 10185.5966796875
$ws.Range("L2").Value = 0.8316212296485901
$ws.Range("M2").Value = 0.8090195059776306
$ws.Range("N2").Value = 0.8276892304420471
$ws.Range("O2").Value = 31308.314453125

$ws.Range("A3").Value = "Ampac-2"
$ws.Range("B3").Value = 28432.502269245717
$ws.Range("C3").Value = 233.8112335205078
$ws.Range("D3").Value = 233.55409240722656
$ws.Range("E3").Value = 235.58250427246094
$ws.Range("F3").Value = 1.158085823059082
$ws.Range("G3").Value = 1.7002322673797607
$ws.Range("H3").Value = 2.4075684547424316
$ws.Range("I3").Value = 243.35299682617188
$ws.Range("J3").Value = 242.78810119628906
$ws.Range("K3").Value = 469.93096923828125
$ws.Range("L3").Value = 0.8985390663146973
$ws.Range("M3").Value = 0.6110633611679077
$ws.Range("N3").Value = 0.8282872438430786
$ws.Range("O3").Value = 956.0720825195312

$ws.Range("A4").Value = "Ledena voda-3"
$ws.Range("B4").Value = 168540.08513597416
$ws.Range("C4").Value = 234.4571533203125
$ws.Range("D4").Value = 235.22039794921875
$ws.Range("E4").Value = 236.38121032714844
$ws.Range("F4").Value = 7.332977294921875
$ws.Range("G4").Value = 8.48817253112793
$ws.Range("H4").Value = 7.296930313110352
$ws.Range("I4").Value = 1263.4212646484375
$ws.Range("J4").Value = 1575.8441162109375
$ws.Range("K4").Value = 1292.468994140625
$ws.Range("L4").Value = 0.7348595857620239
$ws.Range("M4").Value = 0.7892672419548035
$ws.Range("N4").Value = 0.7493193745613098
$ws.Range("O4").Value = 4128.52978515625

$ws.Range("A5").Value = "Hladilnici-4"
$ws.Range("B5").Value = 126064.47047247678
$ws.Range("C5").Value = 235.95840454101562
$ws.Range("D5").Value = 234.39183044433594
$ws.Range("E5").Value = 234.66357421875
$ws.Range("F5").Value = 11.818692207336426
$ws.Range("G5").Value = 10.408903121948242
$ws.Range("H5").Value = 12.592012405395508
$ws.Range("I5").Value = 2542.15771484375
$ws.Range("J5").Value = 1959.582763671875
$ws.Range("K5").Value = 2584.9189453125
$ws.Range("L5").Value = 0.9115859270095825
$ws.Range("M5").Value = 0.8031861186027527
$ws.Range("N5").Value = 0.8747946619987488
$ws.Range("O5").Value = 7090.69677734375

$ws.Range("A6").Value = "Kompresorno-5"
$ws.Range("B6").Value = 26854.9239010611
$ws.Range("C6").Value = 235.74732971191406
$ws.Range("D6").Value = 233.9075164794922
$ws.Range("E6").Value = 233.57058715820312
$ws.Range("F6").Value = 47.047027587890625
$ws.Range("G6").Value = 44.54290008544922
$ws.Range("H6").Value = 45.58487319946289
$ws.Range("I6").Value = 8690.9462890625
$ws.Range("J6").Value = 8164.4853515625
$ws.Range("K6").Value = 8207.5693359375
$ws.Range("L6").Value = 0.7833541035652161
$ws.Range("M6").Value = 0.7836212515830994
$ws.Range("N6").Value = 0.7708601951599121
$ws.Range("O6").Value = 25063

$ws.Range("A7").Value = "Priemno-6"
$ws.Range("B7").Value = 24341.95070399082
$ws.Range("C7").Value = 234.18414306640625
$ws.Range("D7").Value = 235.21852111816406
$ws.Range("E7").Value = 236.48153686523438
$ws.Range("F7").Value = 11.929794311523438
$ws.Range("G7").Value = 4.372923851013184
$ws.Range("H7").Value = 3.900940418243408
$ws.Range("I7").Value = 2652.338134765625
$ws.Range("J7").Value = 588.9971923828125
$ws.Range("K7").Value = 481.669677734375
$ws.Range("L7").Value = 0.950146496295929
$ws.Range("M7").Value = 0.5738595128059387
$ws.Range("N7").Value = 0.5201421976089478
$ws.Range("O7").Value = 3724.3369140625

$ws.Range("A8").Value = "Trafo#1-7"
$ws.Range("B8").Value = 116557.51448604243
$ws.Range("C8").Value = 234.64332580566406
$ws.Range("D8").Value = 234.44541931152344
$ws.Range("E8").Value = 236.64598083496094
$ws.Range("F8").Value = 33.484519958496094
$ws.Range("G8").Value = 42.40336608886719
$ws.Range("H8").Value = 38.169071197509766
$ws.Range("I8").Value = 7236.0830078125
$ws.Range("J8").Value = 8743.78125
$ws.Range("K8").Value = 8364.05078125
$ws.Range("L8").Value = 0.9206428527832031
$ws.Range("M8").Value = 0.8795431852340698
$ws.Range("N8").Value = 0.9259892702102661
$ws.Range("O8").Value = 24343.9140625
